$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) Trim the trailing double space to a single space in the last run of
#    the "...and similarly for the x64 folders.  " paragraph.
# ----------------------------------------------------------------------
$d.Content.Find.Execute("folders.  ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "folders. ", 2) | Out-Null

# Locate that paragraph again (text changed already) so we can anchor the
# new paragraph right after it.
$hostPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*x64 folders. *") {
        $hostPara = $cand
    }
}

# ----------------------------------------------------------------------
# 2) Insert a brand-new paragraph right after it, with no list/bullet
#    formatting (plain "Normal" paragraph).
# ----------------------------------------------------------------------
$hostPara.Range.InsertParagraphAfter() | Out-Null
$newIndex = $hostPara.Index + 1
$newPara = $d.Paragraphs.Item($newIndex)
$newPara.Style = "Normal"

$q1 = [char]0x201C
$q2 = [char]0x201D

$boldText  = "Extra note:"
$restText  = " The SAM registration process can be skipped by either:"
$lineA     = "a) Entering " + $q1 + "09332s" + $q2 + " as the registration code in each version of SAM on your computer, or"
$lineB     = "b) Run the command " + $q1 + "regedit" + $q2 + " from the command line, search for SAMnt, and add a new registry string value named " + $q1 + "developer-registration" + $q2 + " with the value 09332s, which will bypass the registration on every version of SAM installed on your computer."

$pos = $newPara.Range.Start

# "Extra note:" (bold)
$ins = $d.Range($pos, $pos)
$ins.InsertAfter($boldText)
$boldRange = $d.Range($pos, $pos + $boldText.Length)
$boldRange.Bold = 1
$pos = $pos + $boldText.Length

# " The SAM registration process can be skipped by either:" (not bold)
$ins = $d.Range($pos, $pos)
$ins.InsertAfter($restText)
$pos = $pos + $restText.Length

# line break + "a) Entering ... or" (single run together)
$ins = $d.Range($pos, $pos)
$ins.InsertAfter([char]11 + $lineA)
$pos = $pos + 1 + $lineA.Length

# line break + "b) Run the command ... computer." (single run together)
$ins = $d.Range($pos, $pos)
$ins.InsertAfter([char]11 + $lineB)
$pos = $pos + 1 + $lineB.Length

# Bookmark "_GoBack" (zero length) right after the text above -- adding it
# moves Word's single hidden "_GoBack" bookmark here (and removes it from
# its old location near "Build Solution").
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Two trailing runs, each a single space.
$ins = $d.Range($pos, $pos)
$ins.InsertAfter(" ")
$pos = $pos + 1

$ins = $d.Range($pos, $pos)
$ins.InsertAfter(" ")
$pos = $pos + 1

Write-Output "done"
